$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Division_summuary")

function Set-CoreMark($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.Value = $text
    # Re-assert the (already-default) theme color; this nudges the engine into
    # materializing a dedicated font record for these cells (mirrors the
    # distinct-but-visually-identical font the workbook ends up with).
    $r.Font.ThemeColor = 1
}

Set-CoreMark "F2"  "434 [1 core]"
Set-CoreMark "F3"  "2200 [1 core]"
Set-CoreMark "F4"  "114 [1 core, strategy to be updated]"
Set-CoreMark "F7"  "38 [1 core]"
Set-CoreMark "F8"  "1267 [1 core]"
Set-CoreMark "F9"  "0 [1 core]"
Set-CoreMark "F10" "1023 [1 core]"
Set-CoreMark "F11" "978 [1 core]"
Set-CoreMark "F12" "3949 [1 core]"
Set-CoreMark "F15" "1942 [1 core]"
Set-CoreMark "F23" "19201 [1 core]"

# Widen column F to fit the new longer values.
$ws.Range("F1").ColumnWidth = 38.166666666666664

# Row 23 grew slightly taller after the edit.
$ws.Rows.Item(23).RowHeight = 19

# Last selected cell before save.
$ws.Range("G22").Select()
